$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

function Set-CellValue($row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Trim trailing cell-mark/paragraph-mark characters so we only replace
    # the visible text, preserving the run's formatting.
    $r.End = $r.End - 1
    $r.Text = $newText
}

Set-CellValue 1 1 "69×29="
Set-CellValue 1 2 "55×29="
Set-CellValue 1 3 "15×91="
Set-CellValue 1 4 "35×20="
Set-CellValue 1 5 "38×41="

Set-CellValue 5 1 "75×32="
Set-CellValue 5 2 "31×30="
Set-CellValue 5 3 "59×94="
Set-CellValue 5 4 "77×23="
Set-CellValue 5 5 "83×42="

Set-CellValue 10 1 "12×69="
Set-CellValue 10 2 "65×44="
Set-CellValue 10 3 "61×56="
Set-CellValue 10 4 "38×42="
Set-CellValue 10 5 "62×49="

Set-CellValue 15 1 "20×74="
Set-CellValue 15 2 "96×55="
Set-CellValue 15 3 "97×15="
Set-CellValue 15 4 "50×47="
Set-CellValue 15 5 "91×29="

Set-CellValue 20 1 "12×96="
Set-CellValue 20 2 "31×74="
Set-CellValue 20 3 "62×56="
Set-CellValue 20 4 "71×53="
Set-CellValue 20 5 "53×39="
